# chore: update Sheets via scheduled runner
# Refreshes cached market-board figures (currentAveragePrice / price / profit
# columns) for a handful of leve rows across the ALC, ARM, BSM, CRP, CUL,
# GSM, LTW and WVR sheets. Values only - no formulas, formatting or
# structural changes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2680
$ws.Range("J17").Value = 2680
$ws.Range("L17").Value = 8040
$ws.Range("N17").Value = -8376

$ws.Range("H53").Value = 797
$ws.Range("I53").Value = 1037.5
$ws.Range("J53").Value = 41.142857
$ws.Range("K53").Value = 1037.5
$ws.Range("L53").Value = 41.142857
$ws.Range("M53").Value = -400.5
$ws.Range("N53").Value = -1315.142857

$ws.Range("H127").Value = 1158.421
$ws.Range("I127").Value = 556.36365
$ws.Range("J127").Value = 1986.25
$ws.Range("K127").Value = 1669.09095
$ws.Range("L127").Value = 5958.75
$ws.Range("M127").Value = 3290.90905
$ws.Range("N127").Value = -15878.75

$ws.Range("H137").Value = 1100.6833
$ws.Range("I137").Value = 833.4375
$ws.Range("J137").Value = 1406.1072
$ws.Range("K137").Value = 2500.3125
$ws.Range("L137").Value = 4218.321599999999
$ws.Range("M137").Value = 49.6875
$ws.Range("N137").Value = -9318.321599999999

$ws.Range("H138").Value = 1408.05
$ws.Range("I138").Value = 937.5
$ws.Range("J138").Value = 1511.3414
$ws.Range("K138").Value = 2812.5
$ws.Range("L138").Value = 4534.0242
$ws.Range("M138").Value = 2327.5
$ws.Range("N138").Value = -14814.0242

$ws.Range("H141").Value = 586.65216
$ws.Range("I141").Value = 591.0454999999999
$ws.Range("J141").Value = 490
$ws.Range("K141").Value = 1773.1365
$ws.Range("L141").Value = 1470
$ws.Range("M141").Value = 3406.8635
$ws.Range("N141").Value = -11830

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3370.106
$ws.Range("I32").Value = 3014.7144
$ws.Range("J32").Value = 10833.333
$ws.Range("K32").Value = 3014.7144
$ws.Range("L32").Value = 10833.333
$ws.Range("M32").Value = -2727.7144
$ws.Range("N32").Value = -11407.333

$ws.Range("H132").Value = 1358.8302
$ws.Range("I132").Value = 1092.4048
$ws.Range("K132").Value = 3277.2144
$ws.Range("M132").Value = -747.2143999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2344.3333
$ws.Range("J20").Value = 2052.6667
$ws.Range("L20").Value = 2052.6667
$ws.Range("N20").Value = -2546.6667

$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H105").Value = 66668228
$ws.Range("I105").Value = 76924800
$ws.Range("J105").Value = 495
$ws.Range("K105").Value = 76924800
$ws.Range("L105").Value = 495
$ws.Range("M105").Value = -76923053
$ws.Range("N105").Value = -3989

$ws.Range("H107").Value = 1400.4286
$ws.Range("I107").Value = 1518.7273
$ws.Range("K107").Value = 1518.7273
$ws.Range("M107").Value = 401.2727

$ws.Range("H134").Value = 4698.1387
$ws.Range("I134").Value = 1111.862
$ws.Range("J134").Value = 19555.572
$ws.Range("K134").Value = 3335.586
$ws.Range("L134").Value = 58666.716
$ws.Range("M134").Value = -800.5860000000002
$ws.Range("N134").Value = -63736.716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 58778.332
$ws.Range("I22").Value = 566
$ws.Range("K22").Value = 566
$ws.Range("M22").Value = -216

$ws.Range("H31").Value = 2298.1667
$ws.Range("I31").Value = 2374.5293
$ws.Range("K31").Value = 2374.5293
$ws.Range("M31").Value = -2079.5293

$ws.Range("H34").Value = 2298.1667
$ws.Range("I34").Value = 2374.5293
$ws.Range("K34").Value = 2374.5293
$ws.Range("M34").Value = -2172.5293

$ws.Range("H58").Value = 630.6667
$ws.Range("I58").Value = 565.6389
$ws.Range("K58").Value = 565.6389
$ws.Range("M58").Value = -362.6389

$ws.Range("H94").Value = 600.125
$ws.Range("J94").Value = 624.1667
$ws.Range("L94").Value = 624.1667
$ws.Range("N94").Value = -1526.1667

$ws.Range("H132").Value = 2328.5557
$ws.Range("I132").Value = 1857.95
$ws.Range("J132").Value = 3673.1428
$ws.Range("K132").Value = 5573.85
$ws.Range("L132").Value = 11019.4284
$ws.Range("M132").Value = -3043.85
$ws.Range("N132").Value = -16079.4284

$ws.Range("H135").Value = 31999.5
$ws.Range("J135").Value = 31999.5
$ws.Range("L135").Value = 31999.5
$ws.Range("N135").Value = -42139.5

$ws.Range("H136").Value = 630.6667
$ws.Range("I136").Value = 565.6389
$ws.Range("K136").Value = 1696.9167
$ws.Range("M136").Value = 853.0832999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 7146.2
$ws.Range("J107").Value = 10520.7
$ws.Range("L107").Value = 31562.1
$ws.Range("N107").Value = -35402.10000000001

$ws.Range("H122").Value = 839.5
$ws.Range("I122").Value = 535
$ws.Range("J122").Value = 941
$ws.Range("K122").Value = 4815
$ws.Range("L122").Value = 8469
$ws.Range("N122").Value = -13369
$ws.Range("M122").Value = -2365

$ws.Range("H131").Value = 28572726
$ws.Range("I131").Value = 71428780
$ws.Range("J131").Value = 2022.619
$ws.Range("K131").Value = 214286340
$ws.Range("L131").Value = 6067.857
$ws.Range("M131").Value = -214281300
$ws.Range("N131").Value = -16147.857

$ws.Range("H140").Value = 23124.074
$ws.Range("J140").Value = 2914
$ws.Range("L140").Value = 8742
$ws.Range("N140").Value = -19102

$ws.Range("H141").Value = 3010
$ws.Range("I141").Value = 3010
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 9030
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3850
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H130").Value = 34696
$ws.Range("J130").Value = 34696
$ws.Range("L130").Value = 34696
$ws.Range("N130").Value = -44736

$ws.Range("H132").Value = 1972.8125
$ws.Range("I132").Value = 1494.5454
$ws.Range("J132").Value = 3025
$ws.Range("K132").Value = 4483.6362
$ws.Range("L132").Value = 9075
$ws.Range("M132").Value = -1953.6362
$ws.Range("N132").Value = -14135

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1073.4546
$ws.Range("I16").Value = 949.5
$ws.Range("J16").Value = 1404
$ws.Range("K16").Value = 949.5
$ws.Range("L16").Value = 1404
$ws.Range("M16").Value = -779.5
$ws.Range("N16").Value = -1744

$ws.Range("H22").Value = 1017
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1125.5
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 1125.5
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -1715.5

$ws.Range("H27").Value = 1017
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 1125.5
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 1125.5
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -1339.5

$ws.Range("H31").Value = 1873.8889
$ws.Range("I31").Value = 857.5
$ws.Range("J31").Value = 2687
$ws.Range("K31").Value = 857.5
$ws.Range("L31").Value = 2687
$ws.Range("M31").Value = -609.5
$ws.Range("N31").Value = -3183

$ws.Range("H46").Value = 1685.7142
$ws.Range("J46").Value = 1960
$ws.Range("L46").Value = 1960
$ws.Range("N46").Value = -2336

$ws.Range("H121").Value = 40420
$ws.Range("J121").Value = 40420
$ws.Range("L121").Value = 40420
$ws.Range("M121").Value = -43914

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 17800
$ws.Range("J64").Value = 17800
$ws.Range("L64").Value = 17800
$ws.Range("N64").Value = -18296

$ws.Range("H67").Value = 17800
$ws.Range("J67").Value = 17800
$ws.Range("L67").Value = 17800
$ws.Range("N67").Value = -19516

$ws.Range("H113").Value = 533.5
$ws.Range("I113").Value = 390.7143
$ws.Range("K113").Value = 1172.1429
$ws.Range("M113").Value = 997.8571000000002
